# Auto-generated: apply cryptos price/volume updates for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '42.365.30'; E = '  +1.32%  '; ForceText = $false }
    @{ Row = 3; D = '2.286.84'; E = '  +0.47%  '; ForceText = $false }
    @{ Row = 5; D = '156.02'; E = '  +15,479.99%  '; ForceText = $true }
    @{ Row = 6; D = '306.02'; E = '  +0.02%  '; ForceText = $true }
    @{ Row = 7; D = '96.82'; E = '  +4.60%  '; ForceText = $true }
    @{ Row = 8; D = '0.532'; E = '  -0.03%  '; ForceText = $true }
    @{ Row = 9; D = $null; E = '  -0.05%  '; ForceText = $false }
    @{ Row = 10; D = '0.496'; E = '  +1.87%  '; ForceText = $true }
    @{ Row = 11; D = '35.50'; E = '  +8.48%  '; ForceText = $true }
    @{ Row = 12; D = $null; E = '  +1.04%  '; ForceText = $false }
    @{ Row = 13; D = $null; E = '  -1.75%  '; ForceText = $false }
    @{ Row = 14; D = '6.72'; E = '  +0.75%  '; ForceText = $true }
    @{ Row = 15; D = '2.641.87'; E = '  +0.53%  '; ForceText = $false }
    @{ Row = 16; D = '14.53'; E = '  +1.56%  '; ForceText = $true }
    @{ Row = 17; D = '2.277.68'; E = '  -1.05%  '; ForceText = $false }
    @{ Row = 18; D = '0.796'; E = '  +4.15%  '; ForceText = $true }
    @{ Row = 19; D = '42.229.73'; E = '  +1.12%  '; ForceText = $false }
    @{ Row = 20; D = '12.90'; E = '  +5.37%  '; ForceText = $true }
    @{ Row = 21; D = '0.0₃0919'; E = '  +1.13%  '; ForceText = $false }
    @{ Row = 22; D = '6.03'; E = '  +1.46%  '; ForceText = $true }
    @{ Row = 23; D = '68.23'; E = '  +1.49%  '; ForceText = $true }
    @{ Row = 24; D = '244.08'; E = '  +0.54%  '; ForceText = $true }
    @{ Row = 25; D = $null; E = '  +0.22%  '; ForceText = $false }
    @{ Row = 26; D = $null; E = '  +1.19%  '; ForceText = $false }
    @{ Row = 27; D = $null; E = '  -0.19%  '; ForceText = $false }
    @{ Row = 28; D = $null; E = '  -0.25%  '; ForceText = $false }
    @{ Row = 29; D = '36.37'; E = '  +6.61%  '; ForceText = $true }
    @{ Row = 30; D = '9.73'; E = '  +0.98%  '; ForceText = $true }
    @{ Row = 31; D = '2.10'; E = '  +1.41%  '; ForceText = $true }
    @{ Row = 32; D = '161.38'; E = '  +1.34%  '; ForceText = $true }
    @{ Row = 33; D = '5.37'; E = '  +3.54%  '; ForceText = $true }
    @{ Row = 34; D = $null; E = '  -0.09%  '; ForceText = $false }
    @{ Row = 35; D = '0.0755'; E = '  +0.53%  '; ForceText = $true }
    @{ Row = 36; D = '3.09'; E = '  +1.46%  '; ForceText = $true }
    @{ Row = 37; D = $null; E = '  +4.36%  '; ForceText = $false }
    @{ Row = 38; D = '17.41'; E = '  +2.26%  '; ForceText = $true }
    @{ Row = 39; D = $null; E = '  +0.32%  '; ForceText = $false }
    @{ Row = 40; D = $null; E = '  -0.49%  '; ForceText = $false }
    @{ Row = 41; D = $null; E = '  +0.90%  '; ForceText = $false }
    @{ Row = 42; D = '4.24'; E = '  +7.56%  '; ForceText = $true }
    @{ Row = 43; D = '19.98'; E = '  +0.43%  '; ForceText = $true }
    @{ Row = 44; D = '2.012.33'; E = '  -3.02%  '; ForceText = $false }
    @{ Row = 45; D = $null; E = '  +11.11%  '; ForceText = $false }
    @{ Row = 46; D = '0.0286'; E = '  +2.38%  '; ForceText = $true }
    @{ Row = 47; D = '10.28'; E = '  -0.29%  '; ForceText = $true }
    @{ Row = 48; D = '3.00'; E = '  +2.79%  '; ForceText = $true }
    @{ Row = 49; D = '54.11'; E = '  +4.22%  '; ForceText = $true }
    @{ Row = 50; D = $null; E = '  +0.92%  '; ForceText = $false }
    @{ Row = 51; D = '72.92'; E = '  +0.07%  '; ForceText = $true }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($row, 4)
        if ($u.ForceText) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.D
    }
    $ws.Cells.Item($row, 5).Value = $u.E
}

